$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3320.081
$ws.Range("J17").Value = 3566.2646
$ws.Range("L17").Value = 10698.7938
$ws.Range("N17").Value = -11034.7938

$ws.Range("H113").Value = 2006.8
$ws.Range("I113").Value = 1591.619
$ws.Range("J113").Value = 2629.5715
$ws.Range("K113").Value = 1591.619
$ws.Range("L113").Value = 2629.5715
$ws.Range("M113").Value = 1662.381
$ws.Range("N113").Value = -9137.5715

$ws.Range("H137").Value = 2539.2632
$ws.Range("I137").Value = 2396.037
$ws.Range("J137").Value = 2890.818
$ws.Range("K137").Value = 7188.110999999999
$ws.Range("L137").Value = 8672.454000000002
$ws.Range("M137").Value = -4638.110999999999
$ws.Range("N137").Value = -13772.454

$ws.Range("H139").Value = 45482.145
$ws.Range("I139").Value = 40000
$ws.Range("J139").Value = 49593.75
$ws.Range("K139").Value = 40000
$ws.Range("L139").Value = 49593.75
$ws.Range("M139").Value = -34860
$ws.Range("N139").Value = -59873.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H118").Value = 55639.4
$ws.Range("J118").Value = 55639.4
$ws.Range("L118").Value = 55639.4
$ws.Range("N118").Value = -58953.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6186.6553
$ws.Range("I31").Value = 10013.667
$ws.Range("J31").Value = 3485.2354
$ws.Range("K31").Value = 10013.667
$ws.Range("L31").Value = 3485.2354
$ws.Range("M31").Value = -9718.666999999999
$ws.Range("N31").Value = -4075.2354

$ws.Range("H34").Value = 6186.6553
$ws.Range("I34").Value = 10013.667
$ws.Range("J34").Value = 3485.2354
$ws.Range("K34").Value = 10013.667
$ws.Range("L34").Value = 3485.2354
$ws.Range("M34").Value = -9811.666999999999
$ws.Range("N34").Value = -3889.2354

$ws.Range("H135").Value = 43637
$ws.Range("J135").Value = 43637
$ws.Range("L135").Value = 43637
$ws.Range("N135").Value = -53777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 43.925926
$ws.Range("I2").Value = 20
$ws.Range("J2").Value = 48.086956
$ws.Range("K2").Value = 120
$ws.Range("L2").Value = 288.521736
$ws.Range("M2").Value = -7
$ws.Range("N2").Value = -514.521736

$ws.Range("H4").Value = 9804055
$ws.Range("I4").Value = 9804055
$ws.Range("K4").Value = 29412165
$ws.Range("M4").Value = -29412053

$ws.Range("H9").Value = 23914376
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 23914376
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 71743128
$ws.Range("M9").Value = ""
$ws.Range("N9").Value = -71743576

$ws.Range("H19").Value = 500
$ws.Range("I19").Value = 500
$ws.Range("K19").Value = 1500
$ws.Range("M19").Value = -1326

$ws.Range("H22").Value = 200000820
$ws.Range("J22").Value = 2000
$ws.Range("L22").Value = 6000
$ws.Range("N22").Value = -6338

$ws.Range("H27").Value = 200000820
$ws.Range("J27").Value = 2000
$ws.Range("L27").Value = 6000
$ws.Range("N27").Value = -6204

$ws.Range("H33").Value = 117.888885
$ws.Range("I33").Value = 99.333336
$ws.Range("J33").Value = 127.166664
$ws.Range("K33").Value = 596.000016
$ws.Range("L33").Value = 762.999984
$ws.Range("M33").Value = -313.000016
$ws.Range("N33").Value = -1328.999984

$ws.Range("H39").Value = 28852
$ws.Range("J39").Value = 28852
$ws.Range("L39").Value = 86556
$ws.Range("N39").Value = -87144

$ws.Range("H58").Value = 3116
$ws.Range("J58").Value = 3116
$ws.Range("L58").Value = 9348
$ws.Range("N58").Value = -9604

$ws.Range("H64").Value = 83335544
$ws.Range("I64").Value = 200001200
$ws.Range("J64").Value = 2927.7144
$ws.Range("K64").Value = 600003600
$ws.Range("L64").Value = 8783.143199999999
$ws.Range("M64").Value = -600003330
$ws.Range("N64").Value = -9323.143199999999

$ws.Range("H67").Value = 83335544
$ws.Range("I67").Value = 200001200
$ws.Range("J67").Value = 2927.7144
$ws.Range("K67").Value = 600003600
$ws.Range("L67").Value = 8783.143199999999
$ws.Range("M67").Value = -600002664
$ws.Range("N67").Value = -10655.1432

$ws.Range("H70").Value = 2055.8333
$ws.Range("I70").Value = 1083.75
$ws.Range("K70").Value = 3251.25
$ws.Range("M70").Value = -2936.25

$ws.Range("H73").Value = 2055.8333
$ws.Range("I73").Value = 1083.75
$ws.Range("K73").Value = 3251.25
$ws.Range("M73").Value = -2159.25

$ws.Range("H76").Value = 1660
$ws.Range("I76").Value = 1660
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 4980
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -4597
$ws.Range("N76").Value = ""

$ws.Range("H79").Value = 1660
$ws.Range("I79").Value = 1660
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 4980
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -3654
$ws.Range("N79").Value = ""

$ws.Range("H94").Value = 2601.4
$ws.Range("I94").Value = 1012
$ws.Range("J94").Value = 2998.75
$ws.Range("K94").Value = 3036
$ws.Range("L94").Value = 8996.25
$ws.Range("M94").Value = -2360
$ws.Range("N94").Value = -10348.25

$ws.Range("H100").Value = 3600
$ws.Range("J100").Value = 3600
$ws.Range("L100").Value = 10800
$ws.Range("N100").Value = -12422

$ws.Range("H103").Value = 2504.6667
$ws.Range("J103").Value = 4009.3333
$ws.Range("L103").Value = 12027.9999
$ws.Range("N103").Value = -13785.9999

$ws.Range("H106").Value = 2553.6667
$ws.Range("J106").Value = 2553.6667
$ws.Range("L106").Value = 7661.000100000001
$ws.Range("N106").Value = -9553.000100000001

$ws.Range("H109").Value = 2763.25
$ws.Range("J109").Value = 2909.6
$ws.Range("L109").Value = 8728.799999999999
$ws.Range("N109").Value = -10808.8

$ws.Range("H112").Value = 2828.3333
$ws.Range("J112").Value = 3239.5
$ws.Range("L112").Value = 9718.5
$ws.Range("N112").Value = -11934.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5712.558
$ws.Range("I70").Value = 5174.2856
$ws.Range("J70").Value = 5972.4136
$ws.Range("K70").Value = 5174.2856
$ws.Range("L70").Value = 5972.4136
$ws.Range("M70").Value = -4904.2856
$ws.Range("N70").Value = -6512.4136

$ws.Range("H73").Value = 5712.558
$ws.Range("I73").Value = 5174.2856
$ws.Range("J73").Value = 5972.4136
$ws.Range("K73").Value = 5174.2856
$ws.Range("L73").Value = 5972.4136
$ws.Range("M73").Value = -4238.2856
$ws.Range("N73").Value = -7844.4136

$ws.Range("H80").Value = 7341.5
$ws.Range("I80").Value = 9646.091
$ws.Range("J80").Value = 5391.4614
$ws.Range("K80").Value = 9646.091
$ws.Range("L80").Value = 5391.4614
$ws.Range("M80").Value = -8648.091
$ws.Range("N80").Value = -7387.4614

$ws.Range("H83").Value = 7341.5
$ws.Range("I83").Value = 9646.091
$ws.Range("J83").Value = 5391.4614
$ws.Range("K83").Value = 48230.455
$ws.Range("L83").Value = 26957.307
$ws.Range("M83").Value = -43238.455
$ws.Range("N83").Value = -36941.307

$ws.Range("H97").Value = 1528
$ws.Range("I97").Value = 1927.1428
$ws.Range("J97").Value = 1313.0769
$ws.Range("K97").Value = 1927.1428
$ws.Range("L97").Value = 1313.0769
$ws.Range("M97").Value = -1431.1428
$ws.Range("N97").Value = -2305.0769

$ws.Range("H139").Value = 37242.6
$ws.Range("J139").Value = 37242.6
$ws.Range("L139").Value = 37242.6
$ws.Range("N139").Value = -47522.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 782.64703
$ws.Range("I22").Value = 788.75
$ws.Range("J22").Value = 780.7692
$ws.Range("K22").Value = 788.75
$ws.Range("L22").Value = 780.7692
$ws.Range("M22").Value = -493.75
$ws.Range("N22").Value = -1370.7692

$ws.Range("H27").Value = 782.64703
$ws.Range("I27").Value = 788.75
$ws.Range("J27").Value = 780.7692
$ws.Range("K27").Value = 788.75
$ws.Range("L27").Value = 780.7692
$ws.Range("M27").Value = -681.75
$ws.Range("N27").Value = -994.7692

$ws.Range("H35").Value = 28480
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 28480
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 28480
$ws.Range("M35").Value = ""
$ws.Range("N35").Value = -29152

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = ""
$ws.Range("N25").Value = ""

$ws.Range("H138").Value = 51048
$ws.Range("J138").Value = 51048
$ws.Range("L138").Value = 51048
$ws.Range("N138").Value = -61328
